$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert 9 new rows of fresh price data (row 947-955);
# this pushes the existing rows 947-996 down to 956-1005.
$ws.Range("A947:A955").EntireRow.Insert()

# Row 947
$ws.Range("A947").Value = 5
$ws.Range("B947").Value = 'Macroferia Regional de Talca'
$ws.Range("C947").Value = 'Maule'
$ws.Range("E947").Value = 7
$ws.Range("F947").Value = 'Fruta'
$ws.Range("G947").Value = 100106
$ws.Range("H947").Value = 'Oleaginosos'
$ws.Range("I947").Value = 100106002
$ws.Range("J947").Value = 'Palta'
$ws.Range("D947").Value = 44706
$ws.Range("K947").Value = 'Champion'
$ws.Range("L947").Value = 'Primera'
$ws.Range("M947").Value = 180
$ws.Range("N947").Value = 2000
$ws.Range("O947").Value = 2000
$ws.Range("P947").Value = 2000
$ws.Range("Q947").Value = '$/kilo (en caja de 8 kilos )'
$ws.Range("R947").Value = 'Región de O''Higgins'
$ws.Range("S947").Value = 2000
$ws.Range("T947").Value = 1

# Row 948
$ws.Range("A948").Value = 5
$ws.Range("B948").Value = 'Macroferia Regional de Talca'
$ws.Range("C948").Value = 'Maule'
$ws.Range("E948").Value = 7
$ws.Range("F948").Value = 'Fruta'
$ws.Range("G948").Value = 100106
$ws.Range("H948").Value = 'Oleaginosos'
$ws.Range("I948").Value = 100106002
$ws.Range("J948").Value = 'Palta'
$ws.Range("D948").Value = 44706
$ws.Range("K948").Value = 'Champion'
$ws.Range("L948").Value = 'Segunda'
$ws.Range("M948").Value = 200
$ws.Range("N948").Value = 1800
$ws.Range("O948").Value = 1800
$ws.Range("P948").Value = 1800
$ws.Range("Q948").Value = '$/kilo (en caja de 8 kilos )'
$ws.Range("R948").Value = 'Región de O''Higgins'
$ws.Range("S948").Value = 1800
$ws.Range("T948").Value = 1

# Row 949
$ws.Range("A949").Value = 5
$ws.Range("B949").Value = 'Macroferia Regional de Talca'
$ws.Range("C949").Value = 'Maule'
$ws.Range("E949").Value = 7
$ws.Range("F949").Value = 'Fruta'
$ws.Range("G949").Value = 100106
$ws.Range("H949").Value = 'Oleaginosos'
$ws.Range("I949").Value = 100106002
$ws.Range("J949").Value = 'Palta'
$ws.Range("D949").Value = 44706
$ws.Range("K949").Value = 'Champion'
$ws.Range("L949").Value = 'Tercera'
$ws.Range("M949").Value = 150
$ws.Range("N949").Value = 1500
$ws.Range("O949").Value = 1500
$ws.Range("P949").Value = 1500
$ws.Range("Q949").Value = '$/kilo (en caja de 8 kilos )'
$ws.Range("R949").Value = 'Región de O''Higgins'
$ws.Range("S949").Value = 1500
$ws.Range("T949").Value = 1

# Row 950
$ws.Range("A950").Value = 5
$ws.Range("B950").Value = 'Macroferia Regional de Talca'
$ws.Range("C950").Value = 'Maule'
$ws.Range("E950").Value = 7
$ws.Range("F950").Value = 'Fruta'
$ws.Range("G950").Value = 100106
$ws.Range("H950").Value = 'Oleaginosos'
$ws.Range("I950").Value = 100106002
$ws.Range("J950").Value = 'Palta'
$ws.Range("D950").Value = 44706
$ws.Range("K950").Value = 'Hass'
$ws.Range("L950").Value = 'Especial'
$ws.Range("M950").Value = 190
$ws.Range("N950").Value = 3300
$ws.Range("O950").Value = 3300
$ws.Range("P950").Value = 3300
$ws.Range("Q950").Value = '$/kilo (en caja de 17 kilos)'
$ws.Range("R950").Value = 'Cabildo'
$ws.Range("S950").Value = 3300
$ws.Range("T950").Value = 1

# Row 951
$ws.Range("A951").Value = 5
$ws.Range("B951").Value = 'Macroferia Regional de Talca'
$ws.Range("C951").Value = 'Maule'
$ws.Range("E951").Value = 7
$ws.Range("F951").Value = 'Fruta'
$ws.Range("G951").Value = 100106
$ws.Range("H951").Value = 'Oleaginosos'
$ws.Range("I951").Value = 100106002
$ws.Range("J951").Value = 'Palta'
$ws.Range("D951").Value = 44706
$ws.Range("K951").Value = 'Hass'
$ws.Range("L951").Value = 'Primera'
$ws.Range("M951").Value = 200
$ws.Range("N951").Value = 20000
$ws.Range("O951").Value = 20000
$ws.Range("P951").Value = 20000
$ws.Range("Q951").Value = '$/bandeja 10 kilos'
$ws.Range("R951").Value = 'Perú'
$ws.Range("S951").Value = 2000
$ws.Range("T951").Value = 10

# Row 952
$ws.Range("A952").Value = 5
$ws.Range("B952").Value = 'Macroferia Regional de Talca'
$ws.Range("C952").Value = 'Maule'
$ws.Range("E952").Value = 7
$ws.Range("F952").Value = 'Fruta'
$ws.Range("G952").Value = 100106
$ws.Range("H952").Value = 'Oleaginosos'
$ws.Range("I952").Value = 100106002
$ws.Range("J952").Value = 'Palta'
$ws.Range("D952").Value = 44706
$ws.Range("K952").Value = 'Hass'
$ws.Range("L952").Value = 'Segunda'
$ws.Range("M952").Value = 200
$ws.Range("N952").Value = 3100
$ws.Range("O952").Value = 3100
$ws.Range("P952").Value = 3100
$ws.Range("Q952").Value = '$/kilo (en caja de 17 kilos)'
$ws.Range("R952").Value = 'Cabildo'
$ws.Range("S952").Value = 3100
$ws.Range("T952").Value = 1

# Row 953
$ws.Range("A953").Value = 5
$ws.Range("B953").Value = 'Macroferia Regional de Talca'
$ws.Range("C953").Value = 'Maule'
$ws.Range("E953").Value = 7
$ws.Range("F953").Value = 'Fruta'
$ws.Range("G953").Value = 100106
$ws.Range("H953").Value = 'Oleaginosos'
$ws.Range("I953").Value = 100106002
$ws.Range("J953").Value = 'Palta'
$ws.Range("D953").Value = 44706
$ws.Range("K953").Value = 'Hass'
$ws.Range("L953").Value = 'Tercera'
$ws.Range("M953").Value = 250
$ws.Range("N953").Value = 18000
$ws.Range("O953").Value = 18000
$ws.Range("P953").Value = 18000
$ws.Range("Q953").Value = '$/bandeja 10 kilos'
$ws.Range("R953").Value = 'Perú'
$ws.Range("S953").Value = 1800
$ws.Range("T953").Value = 10

# Row 954
$ws.Range("A954").Value = 5
$ws.Range("B954").Value = 'Macroferia Regional de Talca'
$ws.Range("C954").Value = 'Maule'
$ws.Range("E954").Value = 7
$ws.Range("F954").Value = 'Fruta'
$ws.Range("G954").Value = 100106
$ws.Range("H954").Value = 'Oleaginosos'
$ws.Range("I954").Value = 100106002
$ws.Range("J954").Value = 'Palta'
$ws.Range("D954").Value = 44706
$ws.Range("K954").Value = 'Hass'
$ws.Range("L954").Value = 'Segunda'
$ws.Range("M954").Value = 150
$ws.Range("N954").Value = 2800
$ws.Range("O954").Value = 2800
$ws.Range("P954").Value = 2800
$ws.Range("Q954").Value = '$/kilo (en caja de 17 kilos)'
$ws.Range("R954").Value = 'Cabildo'
$ws.Range("S954").Value = 2800
$ws.Range("T954").Value = 1

# Row 955
$ws.Range("A955").Value = 5
$ws.Range("B955").Value = 'Macroferia Regional de Talca'
$ws.Range("C955").Value = 'Maule'
$ws.Range("E955").Value = 7
$ws.Range("F955").Value = 'Fruta'
$ws.Range("G955").Value = 100106
$ws.Range("H955").Value = 'Oleaginosos'
$ws.Range("I955").Value = 100106002
$ws.Range("J955").Value = 'Palta'
$ws.Range("D955").Value = 44706
$ws.Range("K955").Value = 'Hass'
$ws.Range("L955").Value = 'Tercera'
$ws.Range("M955").Value = 180
$ws.Range("N955").Value = 15000
$ws.Range("O955").Value = 15000
$ws.Range("P955").Value = 15000
$ws.Range("Q955").Value = '$/bandeja 10 kilos'
$ws.Range("R955").Value = 'Perú'
$ws.Range("S955").Value = 1500
$ws.Range("T955").Value = 10
